# EBEGU-2029 Namensänderung Jugendamt zu Familie & Quartier Stadt Bern
#
# Renames the report title from "Mitarbeiterinnen Jugendamt" to
# "MitarbeiterInnen Familie & Quartier Stadt Bern" and moves the active
# selection from E8 to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the report title in cell A1.
$ws.Range("A1").Value = "MitarbeiterInnen Familie & Quartier Stadt Bern"

# Move the active selection to A2 (was E8).
[void]$ws.Range("A2").Select()
